$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversion del dia" text block ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText.Replace(
    "1000 Bs = 1.66 = 6079.65 pesos",
    "1000 Bs = 1.85 = 6773.72 pesos"
).Replace(
    "6079.65 pesos = 1.65 = 904.96 Bs",
    "6773.72 pesos = 1.84 = 1011.0 Bs"
)
$cellA1.Value = $newText

# --- tasas sheet: update N10, O10, N12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 541.8
$ws2.Range("O10").Value = 3670
$ws2.Range("N12").Value = 3685
